$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.046010529502491
$ws.Range("D2").Value = 1.050006871171725
$ws.Range("E2").Value = 1.058698708362374
$ws.Range("F2").Value = 1.064156674700211
$ws.Range("I2").Value = 1.037106215224382
$ws.Range("J2").Value = 1.051067206422133
$ws.Range("K2").Value = 1.052762347095428
$ws.Range("L2").Value = 1.061430230530892
$ws.Range("M2").Value = 1.066873371794892
$ws.Range("N2").Value = 1.052559842414346

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.047077892408904
$ws.Range("D3").Value = 1.050812797914315
$ws.Range("E3").Value = 1.059716652479945
$ws.Range("F3").Value = 1.065171196435683
$ws.Range("I3").Value = 1.03727462464952
$ws.Range("J3").Value = 1.051781882566756
$ws.Range("K3").Value = 1.053380650500119
$ws.Range("L3").Value = 1.062261742371102
$ws.Range("M3").Value = 1.067702544301359
$ws.Range("N3").Value = 1.05327553348107

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.047768975451466
$ws.Range("D4").Value = 1.051334436940152
$ws.Range("E4").Value = 1.060376074711297
$ws.Range("F4").Value = 1.06582833956648
$ws.Range("I4").Value = 1.037382278197022
$ws.Range("J4").Value = 1.052244169760625
$ws.Range("K4").Value = 1.053780231084122
$ws.Range("L4").Value = 1.062799918534404
$ws.Range("M4").Value = 1.068239140285755
$ws.Range("N4").Value = 1.053738477175778

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.048059609014307
$ws.Range("D5").Value = 1.051553769133651
$ws.Range("E5").Value = 1.060653473327515
$ws.Range("F5").Value = 1.06610476451776
$ws.Range("I5").Value = 1.037427219709483
$ws.Range("J5").Value = 1.05243847743527
$ws.Range("K5").Value = 1.053948093702892
$ws.Range("L5").Value = 1.063026198900481
$ws.Range("M5").Value = 1.068464740367983
$ws.Range("N5").Value = 1.053933060789609

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.048108413659314
$ws.Range("D6").Value = 1.0515905980092
$ws.Range("E6").Value = 1.060700060165327
$ws.Range("F6").Value = 1.066151186960212
$ws.Range("I6").Value = 1.037434747044046
$ws.Range("J6").Value = 1.052471100314113
$ws.Range("K6").Value = 1.053976271458149
$ws.Range("L6").Value = 1.063064194176657
$ws.Range("M6").Value = 1.0685026204828
$ws.Range("N6").Value = 1.053965729996682

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.047772858509792
$ws.Range("D7").Value = 1.051337367530672
$ws.Range("E7").Value = 1.06037978062854
$ws.Range("F7").Value = 1.065832032534677
$ws.Range("I7").Value = 1.037382879949977
$ws.Range("J7").Value = 1.052246766258679
$ws.Range("K7").Value = 1.053782474549111
$ws.Range("L7").Value = 1.062802941982239
$ws.Range("M7").Value = 1.06824215470646
$ws.Range("N7").Value = 1.053741077361157

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.0463711610907
$ws.Range("D8").Value = 1.050279206183547
$ws.Range("E8").Value = 1.0590425728165
$ws.Range("F8").Value = 1.064499395877375
$ws.Range("I8").Value = 1.037163402788369
$ws.Range("J8").Value = 1.051308766583491
$ws.Range("K8").Value = 1.052971409413098
$ws.Range("L8").Value = 1.061711216080727
$ws.Range("M8").Value = 1.067153580482063
$ws.Range("N8").Value = 1.052801745618841

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.043904475858023
$ws.Range("D9").Value = 1.04841578271728
$ws.Range("E9").Value = 1.056691972271872
$ws.Range("F9").Value = 1.062156358608988
$ws.Range("I9").Value = 1.036766572862163
$ws.Range("J9").Value = 1.049654717975492
$ws.Range("K9").Value = 1.051538383538415
$ws.Range("L9").Value = 1.059788495443386
$ws.Range("M9").Value = 1.065235907166782
$ws.Range("N9").Value = 1.051145348072096

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.042262235318244
$ws.Range("D10").Value = 1.047174353156404
$ws.Range("E10").Value = 1.055128796418796
$ws.Range("F10").Value = 1.060597898455817
$ws.Range("I10").Value = 1.036495261067356
$ws.Range("J10").Value = 1.048551257943479
$ws.Range("K10").Value = 1.050580496480448
$ws.Range("L10").Value = 1.058507413962719
$ws.Range("M10").Value = 1.063957853911985
$ws.Range("N10").Value = 1.050040321000268

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.041551653499132
$ws.Range("D11").Value = 1.0466370130592
$ws.Range("E11").Value = 1.054452854403592
$ws.Range("F11").Value = 1.059923921831505
$ws.Range("I11").Value = 1.036376181219839
$ws.Range("J11").Value = 1.04807327243659
$ws.Range("K11").Value = 1.050165126237377
$ws.Range("L11").Value = 1.05795287163773
$ws.Range("M11").Value = 1.063404543800614
$ws.Range("N11").Value = 1.049561656699134

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.041287789731161
$ws.Range("D12").Value = 1.046437452912384
$ws.Range("E12").Value = 1.054201918707939
$ws.Range("F12").Value = 1.059673704742557
$ws.Range("I12").Value = 1.036331709569444
$ws.Range("J12").Value = 1.047895700652487
$ws.Range("K12").Value = 1.050010749773371
$ws.Range("L12").Value = 1.057746916576181
$ws.Range("M12").Value = 1.063199034631009
$ws.Range("N12").Value = 1.049383832742731

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.04134438587793
$ws.Range("D13").Value = 1.046480257785726
$ws.Range("E13").Value = 1.05425573895594
$ws.Range("F13").Value = 1.059727371362921
$ws.Range("I13").Value = 1.036341259749359
$ws.Range("J13").Value = 1.047933791604641
$ws.Range("K13").Value = 1.050043868089134
$ws.Range("L13").Value = 1.057791093422916
$ws.Range("M13").Value = 1.063243116365644
$ws.Range("N13").Value = 1.049421977788406

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.041529840847652
$ws.Range("D14").Value = 1.046620516700167
$ws.Range("E14").Value = 1.054432109136682
$ws.Range("F14").Value = 1.059903236197541
$ws.Range("I14").Value = 1.03637251007545
$ws.Range("J14").Value = 1.048058594839771
$ws.Range("K14").Value = 1.050152367252762
$ws.Range("L14").Value = 1.057935846781879
$ws.Range("M14").Value = 1.063387556042808
$ws.Range("N14").Value = 1.049546958258445

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.041644116104377
$ws.Range("D15").Value = 1.046706939095494
$ws.Range("E15").Value = 1.054540795095801
$ws.Range("F15").Value = 1.060011609284199
$ws.Range("I15").Value = 1.036391732626983
$ws.Range("J15").Value = 1.048135486706899
$ws.Range("K15").Value = 1.050219205329037
$ws.Range("L15").Value = 1.058025037648417
$ws.Range("M15").Value = 1.063476552077059
$ws.Range("N15").Value = 1.049623959320844

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.042309405228028
$ws.Range("D16").Value = 1.047210019078521
$ws.Range("E16").Value = 1.055173675965118
$ws.Range("F16").Value = 1.0606426459592
$ws.Range("I16").Value = 1.036503130318714
$ws.Range("J16").Value = 1.048582976505261
$ws.Range("K16").Value = 1.050608050662039
$ws.Range("L16").Value = 1.058544220823113
$ws.Range("M16").Value = 1.063994577335027
$ws.Range("N16").Value = 1.050072084606047

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.042726862484902
$ws.Range("D17").Value = 1.047525643854325
$ws.Range("E17").Value = 1.055570913018491
$ws.Range("F17").Value = 1.061038705757654
$ws.Range("I17").Value = 1.036572578977738
$ws.Range("J17").Value = 1.048863627089626
$ws.Range("K17").Value = 1.050851803100185
$ws.Range("L17").Value = 1.058869937720695
$ws.Range("M17").Value = 1.064319546790636
$ws.Range("N17").Value = 1.050353133746433

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042970408335057
$ws.Range("D18").Value = 1.047709762375677
$ws.Range("E18").Value = 1.055802703734494
$ws.Range("F18").Value = 1.061269802411129
$ws.Range("I18").Value = 1.036612932803998
$ws.Range("J18").Value = 1.049027308432386
$ws.Range("K18").Value = 1.050993921961418
$ws.Range("L18").Value = 1.059059939674392
$ws.Range("M18").Value = 1.064509105144583
$ws.Range("N18").Value = 1.050517047535474

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.043053459674494
$ws.Range("D19").Value = 1.047772545390178
$ws.Range("E19").Value = 1.055881753469321
$ws.Range("F19").Value = 1.061348614232962
$ws.Range("I19").Value = 1.036626666217732
$ws.Range("J19").Value = 1.049083116578872
$ws.Range("K19").Value = 1.051042370999144
$ws.Range("L19").Value = 1.059124728266214
$ws.Range("M19").Value = 1.064573741184249
$ws.Range("N19").Value = 1.05057293493593

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.042682068090643
$ws.Range("D20").Value = 1.047491778235843
$ws.Range("E20").Value = 1.05552828402593
$ws.Range("F20").Value = 1.060996203837232
$ws.Range("I20").Value = 1.036565143762196
$ws.Range("J20").Value = 1.048833517734517
$ws.Range("K20").Value = 1.050825656746475
$ws.Range("L20").Value = 1.058834989639347
$ws.Range("M20").Value = 1.064284679691416
$ws.Range("N20").Value = 1.050322981632586

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.041475226829687
$ws.Range("D21").Value = 1.046579213059159
$ws.Range("E21").Value = 1.054380168665715
$ws.Range("F21").Value = 1.059851444853275
$ws.Range("I21").Value = 1.036363314258665
$ws.Range("J21").Value = 1.048021844131362
$ws.Range("K21").Value = 1.050120419418851
$ws.Range("L21").Value = 1.057893219786315
$ws.Range("M21").Value = 1.063345021744119
$ws.Range("N21").Value = 1.049510155359815

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040716888855797
$ws.Range("D22").Value = 1.046005632005117
$ws.Range("E22").Value = 1.053659108978214
$ws.Range("F22").Value = 1.059132429375218
$ws.Range("I22").Value = 1.036235027082854
$ws.Range("J22").Value = 1.047511358514557
$ws.Range("K22").Value = 1.049676491375322
$ws.Range("L22").Value = 1.057301246105808
$ws.Range("M22").Value = 1.062754307623105
$ws.Range("N22").Value = 1.048998944794879

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.04111885542434
$ws.Range("D23").Value = 1.046309680389463
$ws.Range("E23").Value = 1.054041279784355
$ws.Range("F23").Value = 1.059513522746656
$ws.Range("I23").Value = 1.036303166099734
$ws.Range("J23").Value = 1.04778199109406
$ws.Range("K23").Value = 1.049911874899882
$ws.Range("L23").Value = 1.057615047769193
$ws.Range("M23").Value = 1.063067448020285
$ws.Range("N23").Value = 1.049269961703688

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.042702308594329
$ws.Range("D24").Value = 1.047507080589547
$ws.Range("E24").Value = 1.055547545957498
$ws.Range("F24").Value = 1.061015408373214
$ws.Range("I24").Value = 1.036568503893131
$ws.Range("J24").Value = 1.048847122909752
$ws.Range("K24").Value = 1.050837471336873
$ws.Range("L24").Value = 1.058850781121043
$ws.Range("M24").Value = 1.064300434604357
$ws.Range("N24").Value = 1.05033660612873

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.044541783072143
$ws.Range("D25").Value = 1.048897375800259
$ws.Range("E25").Value = 1.057298975822559
$ws.Range("F25").Value = 1.062761464715503
$ws.Range("I25").Value = 1.036870355671472
$ws.Range("J25").Value = 1.0500824653622
$ws.Range("K25").Value = 1.051909304753343
$ws.Range("L25").Value = 1.060285438176797
$ws.Range("M25").Value = 1.065731604660006
$ws.Range("N25").Value = 1.051573702909157

